$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Acelga" (Femacal de La Calera).
# It belongs chronologically before the existing row 148, so insert a new
# row there; Excel shifts rows 148:278 down to 149:279 automatically.
$ws.Rows("148:148").Insert()

$ws.Range("A148").Value = 3
$ws.Range("B148").Value = "Femacal de La Calera"
$ws.Range("C148").Value = "Coquimbo"
$ws.Range("D148").Value = 44589
$ws.Range("E148").Value = 5
$ws.Range("F148").Value = 100112009
$ws.Range("G148").Value = "Acelga"
$ws.Range("H148").Value = "Sin especificar"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 280
$ws.Range("K148").Value = 2300
$ws.Range("L148").Value = 2500
$ws.Range("M148").Value = 2414
$ws.Range("N148").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O148").Value = "Provincia de Quillota"
$ws.Range("P148").Value = 402
$ws.Range("Q148").Value = 6
$ws.Range("R148").Value = "Hortaliza"
